$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2 value (analkanal -> bananpannkaka)
$ws.Range("C2").Value = "bananpannkaka"

# Add new row 7: Dave, Dave, jeff, 2000 (text), User
$ws.Range("A7").Value = "Dave"
$ws.Range("B7").Value = "Dave"
$ws.Range("C7").Value = "jeff"
$ws.Range("D7").Value = "'2000"
$ws.Range("E7").Value = "User"
